# Apply the crypto price/volume update described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.153.71"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +2.23%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.314.03"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +2.04%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.50"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.10%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.06"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +6.71%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.505"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +2.48%  "

# Row 8
$ws.Range("E8").Value = "  -0.06%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.518"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +5.75%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.08"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +8.97%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0796"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.20%  "

# Row 12
$ws.Range("E12").Value = "  +3.51%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.93"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +14.32%  "

# Row 14
$ws.Range("E14").Value = "  +4.03%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.687.15"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.59%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.327.11"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.22%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.815"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +4.31%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.044.72"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.15%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.72"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +9.55%  "

# Row 20
$ws.Range("E20").Value = "  +3.45%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0906"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.98%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.87"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.53%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.81"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.17%  "

# Row 24
$ws.Range("E24").Value = "  +13.47%  "

# Row 25
$ws.Range("E25").Value = "  +0.77%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.17%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.88"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +4.34%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.21"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -3.46%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "167.75"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.31%  "

# Row 30
$ws.Range("E30").Value = "  +2.78%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.28"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.28%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.05%  "

# Row 33
$ws.Range("B33").Value = "RenderToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.75"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +4.76%  "

# Row 34
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.04"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.51%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.30"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +3.72%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.42"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.70%  "

# Row 37
$ws.Range("E37").Value = "  +1.51%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.103"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +4.50%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.84"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.64%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.80"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +5.07%  "

# Row 41
$ws.Range("E41").Value = "  +1.68%  "

# Row 42
$ws.Range("E42").Value = "  -3.97%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.994.77"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.36%  "

# Row 44
$ws.Range("E44").Value = "  +4.54%  "

# Row 45
$ws.Range("E45").Value = "  +8.29%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.66"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.52%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.90"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +4.91%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "56.22"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +7.58%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.521.77"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.15%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.54"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +4.30%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.59"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.51%  "
